# Daily attendance processing - 2026-01-01 14:58:56
# Swap the order of the "Recorded By" names in column G:
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

# Column G = "Recorded By"
$colIndex = 7

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colIndex)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
